$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 225778
$ws.Range("E2").Value = 13108
$ws.Range("F2").Value = 13108
$ws.Range("G2").Value = 11599
$ws.Range("H2").Value = 8540
$ws.Range("I2").Value = 8679
$ws.Range("J2").Value = -139
$ws.Range("K2").Value = 181276
$ws.Range("L2").Value = 58618
$ws.Range("M2").Value = 122659
$ws.Range("N2").Value = 121399
$ws.Range("O2").Value = 1259
$ws.Range("P2").Value = 3695
$ws.Range("Q2").Value = 19944
$ws.Range("R2").Value = -19907
$ws.Range("S2").Value = -4377
$ws.Range("T2").Value = 14113
$ws.Range("U2").Value = 5832
$ws.Range("V2").Value = 29336
$ws.Range("W2").Value = 5.81
$ws.Range("X2").Value = 3.78
$ws.Range("Y2").Value = 7.31
$ws.Range("Z2").Value = 4.8
$ws.Range("AA2").Value = 47.79
$ws.Range("AB2").Value = 3206.6
$ws.Range("AC2").Value = 11745
$ws.Range("AD2").Value = 15.41
$ws.Range("AE2").Value = 165091
$ws.Range("AF2").Value = 1.1
$ws.Range("AG2").Value = 4000
$ws.Range("AH2").Value = 2.21
$ws.Range("AI2").Value = 33.93
$ws.Range("AJ2").Value = 66271100

# Row 3
$ws.Range("D3").Value = 202066
$ws.Range("E3").Value = 18236
$ws.Range("F3").Value = 18236
$ws.Range("G3").Value = 15496
$ws.Range("H3").Value = 11485
$ws.Range("I3").Value = 11530
$ws.Range("J3").Value = -45
$ws.Range("K3").Value = 185787
$ws.Range("L3").Value = 54752
$ws.Range("M3").Value = 131035
$ws.Range("N3").Value = 129915
$ws.Range("O3").Value = 1121
$ws.Range("P3").Value = 3695
$ws.Range("Q3").Value = 31721
$ws.Range("R3").Value = -16978
$ws.Range("S3").Value = -7575
$ws.Range("T3").Value = 16334
$ws.Range("U3").Value = 15387
$ws.Range("V3").Value = 26587
$ws.Range("W3").Value = 9.029999999999999
$ws.Range("X3").Value = 5.68
$ws.Range("Y3").Value = 9.18
$ws.Range("Z3").Value = 6.26
$ws.Range("AA3").Value = 41.78
$ws.Range("AB3").Value = 3434.44
$ws.Range("AC3").Value = 15602
$ws.Range("AD3").Value = 21.06
$ws.Range("AE3").Value = 176671
$ws.Range("AF3").Value = 1.86
$ws.Range("AG3").Value = 4500
$ws.Range("AH3").Value = 1.37
$ws.Range("AI3").Value = 28.73
$ws.Range("AJ3").Value = 66271100

# Row 4
$ws.Range("D4").Value = 206593
$ws.Range("E4").Value = 19919
$ws.Range("F4").Value = 19919
$ws.Range("G4").Value = 16598
$ws.Range("H4").Value = 12810
$ws.Range("I4").Value = 12811
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 204871
$ws.Range("L4").Value = 64361
$ws.Range("M4").Value = 140510
$ws.Range("N4").Value = 139374
$ws.Range("O4").Value = 1136
$ws.Range("P4").Value = 3695
$ws.Range("Q4").Value = 25167
$ws.Range("R4").Value = -17368
$ws.Range("S4").Value = -10073
$ws.Range("T4").Value = 13985
$ws.Range("U4").Value = 11182
$ws.Range("V4").Value = 28906
$ws.Range("W4").Value = 9.640000000000001
$ws.Range("X4").Value = 6.2
$ws.Range("Y4").Value = 9.52
$ws.Range("Z4").Value = 6.56
$ws.Range("AA4").Value = 45.8
$ws.Range("AB4").Value = 3686.21
$ws.Range("AC4").Value = 17336
$ws.Range("AD4").Value = 15.06
$ws.Range("AE4").Value = 192955
$ws.Range("AG4").Value = 5000
$ws.Range("AH4").Value = 1.92
$ws.Range("AI4").Value = 28.73
$ws.Range("AJ4").Value = 66271100

# Row 5
$ws.Range("D5").Value = 256980
$ws.Range("E5").Value = 29285
$ws.Range("F5").Value = 29285
$ws.Range("G5").Value = 25639
$ws.Range("H5").Value = 20220
$ws.Range("I5").Value = 19453
$ws.Range("J5").Value = 767
$ws.Range("K5").Value = 250412
$ws.Range("L5").Value = 87026
$ws.Range("M5").Value = 163386
$ws.Range("N5").Value = 161685
$ws.Range("O5").Value = 1700
$ws.Range("P5").Value = 3914
$ws.Range("Q5").Value = 31807
$ws.Range("R5").Value = -16404
$ws.Range("S5").Value = -7365
$ws.Range("T5").Value = 22526
$ws.Range("U5").Value = 9281
$ws.Range("V5").Value = 30449
$ws.Range("W5").Value = 11.4
$ws.Range("X5").Value = 7.87
$ws.Range("Y5").Value = 12.92
$ws.Range("Z5").Value = 8.880000000000001
$ws.Range("AA5").Value = 53.26
$ws.Range("AB5").Value = 4167.96
$ws.Range("AC5").Value = 24854
$ws.Range("AD5").Value = 16.3
$ws.Range("AE5").Value = 211043
$ws.Range("AF5").Value = 1.92
$ws.Range("AG5").Value = 6000
$ws.Range("AH5").Value = 1.48
$ws.Range("AI5").Value = 23.65
$ws.Range("AJ5").Value = 70592343

# Row 6
$ws.Range("D6").Value = 281830
$ws.Range("E6").Value = 22461
$ws.Range("F6").Value = 22461
$ws.Range("G6").Value = 19400
$ws.Range("H6").Value = 15193
$ws.Range("I6").Value = 14726
$ws.Range("K6").Value = 289441
$ws.Range("L6").Value = 116220
$ws.Range("M6").Value = 173221
$ws.Range("N6").Value = 170830
$ws.Range("P6").Value = 3914
$ws.Range("Q6").Value = 21250
$ws.Range("R6").Value = -36390
$ws.Range("S6").Value = 17938
$ws.Range("T6").Value = 42194
$ws.Range("U6").Value = -20944
$ws.Range("V6").Value = 53211
$ws.Range("W6").Value = 7.97
$ws.Range("X6").Value = 5.39
$ws.Range("Y6").Value = 8.859999999999999
$ws.Range("Z6").Value = 5.63
$ws.Range("AA6").Value = 67.09
$ws.Range("AB6").Value = 4412.01
$ws.Range("AC6").Value = 18812
$ws.Range("AD6").Value = 18.45
$ws.Range("AE6").Value = 222980
$ws.Range("AF6").Value = 1.56
$ws.Range("AG6").Value = 6000
$ws.Range("AH6").Value = 1.73
$ws.Range("AI6").Value = 31.24
$ws.Range("AJ6").Value = 70592343

# Row 7
$ws.Range("D7").Value = 290522
$ws.Range("E7").Value = 10609
$ws.Range("G7").Value = 7807
$ws.Range("H7").Value = 5004
$ws.Range("I7").Value = 4574
$ws.Range("K7").Value = 333953
$ws.Range("L7").Value = 159164
$ws.Range("M7").Value = 174789
$ws.Range("N7").Value = 171981
$ws.Range("P7").Value = 3911
$ws.Range("Q7").Value = 28218
$ws.Range("R7").Value = -57385
$ws.Range("S7").Value = 32066
$ws.Range("T7").Value = 56126
$ws.Range("U7").Value = -35132
$ws.Range("W7").Value = 3.65
$ws.Range("X7").Value = 1.72
$ws.Range("Y7").Value = 2.67
$ws.Range("Z7").Value = 1.6
$ws.Range("AA7").Value = 91.06
$ws.Range("AC7").Value = 5843
$ws.Range("AD7").Value = 57.76
$ws.Range("AE7").Value = 224482
$ws.Range("AF7").Value = 1.5
$ws.Range("AG7").Value = 4033
$ws.Range("AH7").Value = 1.2
$ws.Range("AI7").Value = 62.24

# Row 8
$ws.Range("D8").Value = 352742
$ws.Range("E8").Value = 17227
$ws.Range("G8").Value = 14168
$ws.Range("H8").Value = 10633
$ws.Range("I8").Value = 10129
$ws.Range("K8").Value = 369763
$ws.Range("L8").Value = 187866
$ws.Range("M8").Value = 181897
$ws.Range("N8").Value = 178595
$ws.Range("P8").Value = 3911
$ws.Range("Q8").Value = 30082
$ws.Range("R8").Value = -47486
$ws.Range("S8").Value = 16856
$ws.Range("T8").Value = 47038
$ws.Range("U8").Value = -18075
$ws.Range("W8").Value = 4.88
$ws.Range("X8").Value = 3.02
$ws.Range("Y8").Value = 5.78
$ws.Range("Z8").Value = 3.02
$ws.Range("AA8").Value = 103.28
$ws.Range("AC8").Value = 12939
$ws.Range("AD8").Value = 26.08
$ws.Range("AE8").Value = 233115
$ws.Range("AF8").Value = 1.45
$ws.Range("AG8").Value = 5094
$ws.Range("AH8").Value = 1.51
$ws.Range("AI8").Value = 35.51

# Row 9
$ws.Range("D9").Value = 396734
$ws.Range("E9").Value = 23614
$ws.Range("G9").Value = 20213
$ws.Range("H9").Value = 15281
$ws.Range("I9").Value = 14530
$ws.Range("K9").Value = 395459
$ws.Range("L9").Value = 203443
$ws.Range("M9").Value = 192016
$ws.Range("N9").Value = 188005
$ws.Range("P9").Value = 3911
$ws.Range("Q9").Value = 38980
$ws.Range("R9").Value = -39493
$ws.Range("S9").Value = 6829
$ws.Range("T9").Value = 40022
$ws.Range("U9").Value = -2393
$ws.Range("W9").Value = 5.95
$ws.Range("X9").Value = 3.85
$ws.Range("Y9").Value = 7.93
$ws.Range("Z9").Value = 3.99
$ws.Range("AA9").Value = 105.95
$ws.Range("AC9").Value = 18561
$ws.Range("AD9").Value = 18.18
$ws.Range("AE9").Value = 245398
$ws.Range("AF9").Value = 1.38
$ws.Range("AG9").Value = 6147
$ws.Range("AH9").Value = 1.82
$ws.Range("AI9").Value = 29.87

Write-Output "LG Chem IFRS data updated"